# Updates the cryptocurrency price/volume table on Sheet1 to reflect the
# latest GitHub Actions data refresh (prices, % changes, and a few
# re-sorted coin rows where ranking order changed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, B (Coin), C (Link), D (Price), E (Volume/1h).
# $null means "leave this cell unchanged".
$updates = @(
    @(2, $null, $null, "44.530.60", "  +1.35%  "),
    @(3, $null, $null, "2.249.48", "  +0.82%  "),
    @(4, $null, $null, $null, "  +0.43%  "),
    @(5, $null, $null, "309.02", "  +1.84%  "),
    @(6, $null, $null, "94.98", "  +0.68%  "),
    @(7, $null, $null, $null, "  +1.02%  "),
    @(8, $null, $null, "1.01", "  +0.17%  "),
    @(9, $null, $null, "0.528", "  +2.16%  "),
    @(10, $null, $null, "35.09", "  +2.41%  "),
    @(11, $null, $null, "0.0808", "  +0.85%  "),
    @(12, $null, $null, "7.29", "  +2.56%  "),
    @(13, $null, $null, "0.105", "  +1.39%  "),
    @(14, $null, $null, "2.310.82", "  +2.14%  "),
    @(15, $null, $null, "0.839", "  +3.65%  "),
    @(16, $null, $null, "13.66", "  +2.51%  "),
    @(17, $null, $null, "44.220.86", "  +1.05%  "),
    @(18, $null, $null, "0.0₃0966", "  +1.80%  "),
    @(19, $null, $null, $null, "  +4.78%  "),
    @(20, $null, $null, "12.16", "  +1.51%  "),
    @(21, $null, $null, "65.92", "  +2.33%  "),
    @(22, $null, $null, "238.62", "  +1.23%  "),
    @(23, $null, $null, $null, "  +4.00%  "),
    @(24, $null, $null, "2.01", "  +4.33%  "),
    @(25, $null, $null, $null, "  +0.00%  "),
    @(26, $null, $null, "2.23", "  +5.54%  "),
    @(27, $null, $null, "9.86", "  +0.97%  "),
    @(28, $null, $null, "37.78", "  +5.43%  "),
    @(29, $null, $null, "6.01", "  +2.84%  "),
    @(30, $null, $null, "20.12", "  +1.29%  "),
    @(31, $null, $null, "152.54", "  +0.08%  "),
    @(32, $null, $null, $null, "  +0.12%  "),
    @(33, $null, $null, $null, "  -0.13%  "),
    @(34, $null, $null, "3.18", "  -1.70%  "),
    @(35, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.109", "  +2.68%  "),
    @(36, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.120", "  +2.38%  "),
    @(37, $null, $null, "1.80", "  +2.49%  "),
    @(38, $null, $null, $null, "  +4.75%  "),
    @(39, $null, $null, "14.49", "  -0.47%  "),
    @(40, $null, $null, "3.82", "  +1.02%  "),
    @(41, $null, $null, $null, "  +2.80%  "),
    @(42, $null, $null, $null, "  +0.23%  "),
    @(43, $null, $null, "1.754.07", "  +1.67%  "),
    @(44, "BitcoinSV", "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv", "80.98", "  -3.14%  "),
    @(45, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.193", "  +5.08%  "),
    @(46, $null, $null, "100.02", "  +1.14%  "),
    @(47, $null, $null, "71.05", "  +4.79%  "),
    @(48, $null, $null, "55.64", "  +4.34%  "),
    @(49, "THORChain", "https://coinranking.com/coin/ybmU-kKU+thorchain-rune", "4.88", "  +0.19%  "),
    @(50, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "8.18", "  +3.15%  "),
    @(51, $null, $null, $null, "  +5.72%  ")
)

foreach ($u in $updates) {
    $row = $u[0]
    $coin  = $u[1]
    $link  = $u[2]
    $price = $u[3]
    $vol   = $u[4]

    if ($null -ne $coin) {
        $ws.Cells.Item($row, 2).Value = $coin
    }
    if ($null -ne $link) {
        $ws.Cells.Item($row, 3).Value = $link
    }
    if ($null -ne $price) {
        # Column D holds price text that often looks numeric (e.g. "1.01",
        # "44.530.60"). Force it to remain plain text so formatting such as
        # trailing zeros and thousand-separator dots is preserved exactly,
        # then restore the cell's original (default) style so no visible
        # formatting changes are introduced.
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $price
        $cell.Style = "Normal"
    }
    if ($null -ne $vol) {
        $ws.Cells.Item($row, 5).Value = $vol
    }
}
